# Applies the data updates described by the commit:
#   "Added unit tests for utils and distractor analysis."
# which re-ran the underlying distractor/missing-response analysis and
# refreshed the numbers in the per-category sheets (OM/NV/NR/ND/ALL) and
# the derived "summary" sheet of mv_person_booklet1.xlsx.

$wb = $excel.ActiveWorkbook

function Set-NumberCell($ws, $addr, $val) {
    # Plain numeric write.
    $ws.Range($addr).Value = $val
}

function Set-TextCell($ws, $addr, $val) {
    # Force a genuine text/shared-string cell (not Excel's automatic
    # "numeric-looking string becomes a number" coercion), then drop back
    # to the default style so we don't leave a stray numFmt on the cell.
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# ---------------------------------------------------------------------
# Sheet "OM"
# ---------------------------------------------------------------------
$wsOM = $wb.Worksheets.Item("OM")
Set-NumberCell $wsOM "B2" 81.47
Set-NumberCell $wsOM "B3" 17.06
Set-NumberCell $wsOM "B4" 1.32

# ---------------------------------------------------------------------
# Sheet "NV"
# ---------------------------------------------------------------------
$wsNV = $wb.Worksheets.Item("NV")
Set-NumberCell $wsNV "B2" 94.56
Set-NumberCell $wsNV "B3" 5.29
Set-TextCell   $wsNV "A4" "2"
Set-NumberCell $wsNV "B4" 0.15

# ---------------------------------------------------------------------
# Sheet "NR"
# ---------------------------------------------------------------------
$wsNR = $wb.Worksheets.Item("NR")
Set-NumberCell $wsNR "B2" 53.97
Set-NumberCell $wsNR "B3" 34.56
Set-NumberCell $wsNR "B4" 9.56
Set-NumberCell $wsNR "B5" 1.32
Set-NumberCell $wsNR "B7" 0.29

# ---------------------------------------------------------------------
# Sheet "ND" -- unchanged
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# Sheet "ALL"
# ---------------------------------------------------------------------
$wsALL = $wb.Worksheets.Item("ALL")
Set-NumberCell $wsALL "B2" 39.12
Set-NumberCell $wsALL "B3" 41.47
Set-NumberCell $wsALL "B4" 15
Set-NumberCell $wsALL "B5" 3.38
Set-NumberCell $wsALL "B6" 0.74
Set-NumberCell $wsALL "B7" 0.29

# ---------------------------------------------------------------------
# Sheet "summary" -- derived percentages + category labels
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("summary")

Set-TextCell $wsSummary "B2" "0.06"
Set-TextCell $wsSummary "C2" "0.6"
Set-TextCell $wsSummary "E2" "0.86"
Set-TextCell $wsSummary "A3" "0.44"
Set-TextCell $wsSummary "B3" "0.24"
Set-TextCell $wsSummary "C3" "0.78"
Set-TextCell $wsSummary "E3" "0.88"

Set-TextCell $wsSummary "C4" "0"
Set-TextCell $wsSummary "B6" "2"
